$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-360 all held the date serial 45178
# (2023-09-09) and were bumped forward by one day to 45179 (2023-09-10).
$ws.Range("C2:C360").Value = 45179
